# Updated symbol list on Sat Jan  7 23:24:33 UTC 2023 with GitHub Actions
#
# Refreshes the Price (D) and Volume(1h) (E) columns for the rows whose
# scraped values moved since the last run. The sheet stores these as plain
# text (e.g. "261.35", "0.66%"), so each target cell is pre-formatted as
# Text ("@") before the write to stop Excel's automatic number/percentage
# coercion from turning the literal string into a numeric value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2;  D = "261.35";     E = "0.66%" }
    @{ Row = 3;  D = "27.16";      E = "1.01%" }
    @{ Row = 4;  D = "4.707";      E = "0.71%" }
    @{ Row = 5;  D = "0.06200";    E = "2.18%" }
    @{ Row = 6;  D = "6.731";      E = "0.54%" }
    @{ Row = 7;  D = $null;        E = "-1.17%" }
    @{ Row = 8;  D = "0.9069";     E = "-1.56%" }
    @{ Row = 9;  D = "0.1406";     E = "0.41%" }
    @{ Row = 10; D = "0.04758";    E = "-11.02%" }
    @{ Row = 11; D = "0.07098";    E = "-0.11%" }
    @{ Row = 12; D = "0.03162";    E = "1.18%" }
    @{ Row = 13; D = "0.09057";    E = "-0.96%" }
    @{ Row = 14; D = "0.001537";   E = "0.58%" }
    @{ Row = 15; D = "0.0006142";  E = "1.48%" }
    @{ Row = 16; D = "0.006112";   E = "0.31%" }
    @{ Row = 17; D = "3.467";      E = "-0.48%" }
    @{ Row = 18; D = "3.171";      E = "0.03%" }
    @{ Row = 19; D = "2.178";      E = $null }
    @{ Row = 20; D = $null;        E = "-0.69%" }
    @{ Row = 21; D = $null;        E = "-1.24%" }
    @{ Row = 22; D = "4.124";      E = "0.44%" }
    @{ Row = 23; D = "0.04220";    E = "-0.26%" }
    @{ Row = 24; D = "0.001219";   E = "0.24%" }
    @{ Row = 25; D = "0.004118";   E = "2.33%" }
    @{ Row = 26; D = $null;        E = "0.16%" }
    @{ Row = 40; D = "0.03902";    E = "0.78%" }
    @{ Row = 41; D = "0.1113";     E = "-0.31%" }
    @{ Row = 42; D = "0.004131";   E = "0.57%" }
    @{ Row = 43; D = $null;        E = "-0.67%" }
    @{ Row = 44; D = "0.01345";    E = "-9.99%" }
    @{ Row = 45; D = "0.00005175"; E = "-4.72%" }
    @{ Row = 46; D = $null;        E = "0.16%" }
    @{ Row = 47; D = "0.03592";    E = "-34.12%" }
    @{ Row = 48; D = "0.1698";     E = "28.53%" }
    @{ Row = 49; D = $null;        E = "0.16%" }
    @{ Row = 50; D = $null;        E = "0.16%" }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $cell = $ws.Range("D$($u.Row)")
        $cell.NumberFormat = "@"
        $cell.Value = $u.D
    }
    if ($null -ne $u.E) {
        $cell = $ws.Range("E$($u.Row)")
        $cell.NumberFormat = "@"
        $cell.Value = $u.E
    }
}
